$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read existing values from row 2 (Job_Title, Job_Description) to reuse on the new row
$jobTitle = $ws.Cells.Item(2, 2).Value2
$jobDescription = $ws.Cells.Item(2, 3).Value2

# Add a new job posting row (row 3) for Job_Id = JD_002
$ws.Cells.Item(3, 1).Value2 = "JD_002"
$ws.Cells.Item(3, 2).Value2 = $jobTitle
$ws.Cells.Item(3, 3).Value2 = $jobDescription
$ws.Cells.Item(3, 4).Value2 = 1
$ws.Cells.Item(3, 5).Value2 = 2

# Keep the new row's height consistent with the other data rows (avoid autofit growing it)
$ws.Rows.Item(3).AutoFit()
